$wb = $excel.ActiveWorkbook

# Sheet ALC, row 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 10).Value = 177
$ws.Cells.Item(9, 11).Value = 177
$ws.Cells.Item(9, 12).Value = 177
$ws.Cells.Item(9, 13).Value = -8
$ws.Cells.Item(9, 14).Value = -515

# Sheet ALC, row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 5416.5
$ws.Cells.Item(18, 9).Value = 5299.8
$ws.Cells.Item(18, 10).Value = 6000
$ws.Cells.Item(18, 11).Value = 5299.8
$ws.Cells.Item(18, 12).Value = 6000
$ws.Cells.Item(18, 13).Value = -5015.8

# Sheet ALC, row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 178.58333
$ws.Cells.Item(41, 9).Value = 178.58333
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 11).Value = 178.58333
$ws.Cells.Item(41, 12).Value = 0
$ws.Cells.Item(41, 13).Value = 261.41667

# Sheet ALC, row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 6199
$ws.Cells.Item(43, 9).Value = 8665
$ws.Cells.Item(43, 10).Value = 2500
$ws.Cells.Item(43, 11).Value = 8665
$ws.Cells.Item(43, 12).Value = 2500
$ws.Cells.Item(43, 13).Value = -8596
$ws.Cells.Item(43, 14).Value = -2638

# Sheet ALC, row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 3333.3333
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 3333.3333
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 3333.3333
$ws.Cells.Item(86, 14).Value = -5579.3333
$ws.Cells.Item(86, 13).ClearContents()

# Sheet ALC, row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 3333.3333
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 3333.3333
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 12).Value = 16666.6665
$ws.Cells.Item(89, 14).Value = -27898.6665
$ws.Cells.Item(89, 13).ClearContents()

# Sheet ALC, row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(96, 8).Value = 1156.25
$ws.Cells.Item(96, 9).Value = 1386.6666
$ws.Cells.Item(96, 10).Value = 465
$ws.Cells.Item(96, 11).Value = 4159.9998
$ws.Cells.Item(96, 12).Value = 1395
$ws.Cells.Item(96, 13).Value = -2786.9998

# Sheet ALC, row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 1406.9474
$ws.Cells.Item(98, 9).Value = 1553.1538
$ws.Cells.Item(98, 10).Value = 1090.1666
$ws.Cells.Item(98, 11).Value = 1553.1538
$ws.Cells.Item(98, 12).Value = 1090.1666
$ws.Cells.Item(98, 13).Value = -55.15380000000005
$ws.Cells.Item(98, 14).Value = -4086.1666

# Sheet ALC, row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106, 8).Value = 29879.53
$ws.Cells.Item(106, 9).Value = 29730.133
$ws.Cells.Item(106, 10).Value = 31000
$ws.Cells.Item(106, 11).Value = 29730.133
$ws.Cells.Item(106, 12).Value = 31000
$ws.Cells.Item(106, 13).Value = -29099.133

# Sheet ALC, row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 1406.9474
$ws.Cells.Item(122, 9).Value = 1553.1538
$ws.Cells.Item(122, 10).Value = 1090.1666
$ws.Cells.Item(122, 11).Value = 4659.4614
$ws.Cells.Item(122, 12).Value = 3270.4998
$ws.Cells.Item(122, 13).Value = -2209.4614
$ws.Cells.Item(122, 14).Value = -8170.4998

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 3760.0977
$ws.Cells.Item(138, 9).Value = 3065.7666
$ws.Cells.Item(138, 10).Value = 5653.727
$ws.Cells.Item(138, 11).Value = 9197.299800000001
$ws.Cells.Item(138, 12).Value = 16961.181
$ws.Cells.Item(138, 13).Value = -4057.299800000001

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6374.2607
$ws.Cells.Item(32, 9).Value = 4335.55
$ws.Cells.Item(32, 10).Value = 19965.666
$ws.Cells.Item(32, 11).Value = 4335.55
$ws.Cells.Item(32, 12).Value = 19965.666
$ws.Cells.Item(32, 13).Value = -4048.55

# Sheet ARM, row 36
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(36, 8).Value = 15008.667
$ws.Cells.Item(36, 9).Value = 12513
$ws.Cells.Item(36, 10).Value = 20000
$ws.Cells.Item(36, 11).Value = 12513
$ws.Cells.Item(36, 12).Value = 20000
$ws.Cells.Item(36, 13).Value = -12167

# Sheet ARM, row 62
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(62, 8).Value = 31362.25
$ws.Cells.Item(62, 9).Value = 35000
$ws.Cells.Item(62, 10).Value = 30149.666
$ws.Cells.Item(62, 11).Value = 35000
$ws.Cells.Item(62, 12).Value = 30149.666
$ws.Cells.Item(62, 13).Value = -34376
$ws.Cells.Item(62, 14).Value = -31397.666

# Sheet ARM, row 65
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(65, 8).Value = 31362.25
$ws.Cells.Item(65, 9).Value = 35000
$ws.Cells.Item(65, 10).Value = 30149.666
$ws.Cells.Item(65, 11).Value = 105000
$ws.Cells.Item(65, 12).Value = 90448.99800000001
$ws.Cells.Item(65, 13).Value = -101880
$ws.Cells.Item(65, 14).Value = -96688.99800000001

# Sheet ARM, row 96
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(96, 8).Value = 9583
$ws.Cells.Item(96, 9).Value = 0
$ws.Cells.Item(96, 10).Value = 9583
$ws.Cells.Item(96, 11).Value = 0
$ws.Cells.Item(96, 12).Value = 9583
$ws.Cells.Item(96, 14).Value = -15075

# Sheet ARM, row 101
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(101, 8).Value = 26401
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 10).Value = 26401
$ws.Cells.Item(101, 11).Value = 0
$ws.Cells.Item(101, 12).Value = 26401
$ws.Cells.Item(101, 14).Value = -32891

# Sheet ARM, row 103
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(103, 8).Value = 0
$ws.Cells.Item(103, 9).Value = 0
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 11).Value = 0
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 14).ClearContents()

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1851.3334
$ws.Cells.Item(132, 9).Value = 1922.2
$ws.Cells.Item(132, 10).Value = 1497
$ws.Cells.Item(132, 11).Value = 5766.6
$ws.Cells.Item(132, 12).Value = 4491
$ws.Cells.Item(132, 13).Value = -3236.6
$ws.Cells.Item(132, 14).Value = -9551

# Sheet BSM, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1668.9
$ws.Cells.Item(86, 9).Value = 1562.7142
$ws.Cells.Item(86, 10).Value = 1916.6666
$ws.Cells.Item(86, 11).Value = 1562.7142
$ws.Cells.Item(86, 12).Value = 1916.6666
$ws.Cells.Item(86, 13).Value = -439.7141999999999

# Sheet BSM, row 88
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(88, 8).Value = 26760.875
$ws.Cells.Item(88, 9).Value = 8408
$ws.Cells.Item(88, 10).Value = 29382.715
$ws.Cells.Item(88, 11).Value = 8408
$ws.Cells.Item(88, 12).Value = 29382.715
$ws.Cells.Item(88, 13).Value = -8002
$ws.Cells.Item(88, 14).Value = -30194.715

# Sheet BSM, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 1668.9
$ws.Cells.Item(89, 9).Value = 1562.7142
$ws.Cells.Item(89, 10).Value = 1916.6666
$ws.Cells.Item(89, 11).Value = 7813.571
$ws.Cells.Item(89, 12).Value = 9583.333000000001
$ws.Cells.Item(89, 13).Value = -2197.571

# Sheet BSM, row 91
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(91, 8).Value = 26760.875
$ws.Cells.Item(91, 9).Value = 8408
$ws.Cells.Item(91, 10).Value = 29382.715
$ws.Cells.Item(91, 11).Value = 8408
$ws.Cells.Item(91, 12).Value = 29382.715
$ws.Cells.Item(91, 13).Value = -7004
$ws.Cells.Item(91, 14).Value = -32190.715

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2910.8125
$ws.Cells.Item(134, 9).Value = 2085.889
$ws.Cells.Item(134, 10).Value = 3971.4285
$ws.Cells.Item(134, 11).Value = 6257.667
$ws.Cells.Item(134, 12).Value = 11914.2855
$ws.Cells.Item(134, 13).Value = -3722.667

# Sheet CRP, row 88
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(88, 8).Value = 40097.8
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 40097.8
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 40097.8
$ws.Cells.Item(88, 14).Value = -40909.8

# Sheet CRP, row 91
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(91, 8).Value = 40097.8
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 40097.8
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 12).Value = 40097.8
$ws.Cells.Item(91, 14).Value = -42905.8

# Sheet CRP, row 95
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(95, 8).Value = 34205.875
$ws.Cells.Item(95, 9).Value = 0
$ws.Cells.Item(95, 10).Value = 34205.875
$ws.Cells.Item(95, 11).Value = 0
$ws.Cells.Item(95, 12).Value = 34205.875
$ws.Cells.Item(95, 14).Value = -39697.875

# Sheet CRP, row 103
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(103, 8).Value = 32797.6
$ws.Cells.Item(103, 9).Value = 24997.25
$ws.Cells.Item(103, 10).Value = 63999
$ws.Cells.Item(103, 11).Value = 24997.25
$ws.Cells.Item(103, 12).Value = 63999
$ws.Cells.Item(103, 13).Value = -23825.25

# Sheet CUL, row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 6875407
$ws.Cells.Item(4, 9).Value = 11579207
$ws.Cells.Item(4, 10).Value = 621.8461
$ws.Cells.Item(4, 11).Value = 34737621
$ws.Cells.Item(4, 12).Value = 1865.5383
$ws.Cells.Item(4, 13).Value = -34737509
$ws.Cells.Item(4, 14).Value = -2089.5383

# Sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 14).ClearContents()

# Sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 14).ClearContents()

# Sheet GSM, row 92
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(92, 8).Value = 16812.5
$ws.Cells.Item(92, 9).Value = 25000
$ws.Cells.Item(92, 10).Value = 14083.333
$ws.Cells.Item(92, 11).Value = 25000
$ws.Cells.Item(92, 12).Value = 14083.333
$ws.Cells.Item(92, 13).Value = -23128
$ws.Cells.Item(92, 14).Value = -17827.333

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 36499.31
$ws.Cells.Item(122, 9).Value = 1370.1305
$ws.Cells.Item(122, 10).Value = 171161.17
$ws.Cells.Item(122, 11).Value = 4110.3915
$ws.Cells.Item(122, 12).Value = 513483.51
$ws.Cells.Item(122, 13).Value = -1660.3915

# Sheet GSM, row 134
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(134, 8).Value = 124999.664
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 10).Value = 124999.664
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 12).Value = 374998.992
$ws.Cells.Item(134, 14).Value = -380068.992

# Sheet LTW, row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1543.25
$ws.Cells.Item(16, 9).Value = 1361.6
$ws.Cells.Item(16, 10).Value = 1846
$ws.Cells.Item(16, 11).Value = 1361.6
$ws.Cells.Item(16, 12).Value = 1846
$ws.Cells.Item(16, 13).Value = -1191.6
$ws.Cells.Item(16, 14).Value = -2186

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2288.0908
$ws.Cells.Item(40, 9).Value = 2355.5557
$ws.Cells.Item(40, 10).Value = 1984.5
$ws.Cells.Item(40, 11).Value = 2355.5557
$ws.Cells.Item(40, 12).Value = 1984.5
$ws.Cells.Item(40, 13).Value = -2219.5557

# Sheet LTW, row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 1585
$ws.Cells.Item(82, 9).Value = 876
$ws.Cells.Item(82, 10).Value = 3003
$ws.Cells.Item(82, 11).Value = 876
$ws.Cells.Item(82, 12).Value = 3003
$ws.Cells.Item(82, 13).Value = -515
$ws.Cells.Item(82, 14).Value = -3725

# Sheet LTW, row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 1585
$ws.Cells.Item(85, 9).Value = 876
$ws.Cells.Item(85, 10).Value = 3003
$ws.Cells.Item(85, 11).Value = 876
$ws.Cells.Item(85, 12).Value = 3003
$ws.Cells.Item(85, 13).Value = 372
$ws.Cells.Item(85, 14).Value = -5499

# Sheet LTW, row 87
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(87, 8).Value = 15000
$ws.Cells.Item(87, 9).Value = 15000
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 15000
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 13).Value = -13877

# Sheet LTW, row 90
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(90, 8).Value = 15000
$ws.Cells.Item(90, 9).Value = 15000
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 11).Value = 45000
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 13).Value = -39384

# Sheet WVR, row 3
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3, 8).Value = 2331.3333
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 2331.3333
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 2331.3333
$ws.Cells.Item(3, 14).Value = -2559.3333
$ws.Cells.Item(3, 13).ClearContents()

# Sheet WVR, row 14
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 2832.6667
$ws.Cells.Item(14, 9).Value = 2832.6667
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 2832.6667
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = -2664.6667

# Sheet WVR, row 70
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(70, 8).Value = 54999.5
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 54999.5
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 54999.5
$ws.Cells.Item(70, 14).Value = -55629.5

# Sheet WVR, row 73
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(73, 8).Value = 54999.5
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 54999.5
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 54999.5
$ws.Cells.Item(73, 14).Value = -57183.5

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 4258.647
$ws.Cells.Item(132, 9).Value = 3498.8
$ws.Cells.Item(132, 10).Value = 5344.143
$ws.Cells.Item(132, 11).Value = 10496.4
$ws.Cells.Item(132, 12).Value = 16032.429
$ws.Cells.Item(132, 13).Value = -7966.400000000001
